# Update the "想去人数" (F column) values for matching rows on both the
# "展览" sheet and the "全部类型" sheet, per the latest data refresh
# (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# "展览" sheet: row number -> new F-column value.
$exhibitionRows = @{
    2  = 1887
    6  = 2635
    7  = 175
    10 = 1555
    11 = 539
    13 = 337
    17 = 214
    21 = 190
    23 = 1690
    26 = 23
    28 = 212
}

# "全部类型" sheet: same underlying events, shifted by a couple of rows.
$allTypesRows = @{
    2  = 1887
    7  = 2635
    8  = 175
    11 = 1555
    12 = 539
    14 = 337
    18 = 214
    22 = 190
    24 = 1690
    27 = 23
    29 = 212
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($r in $exhibitionRows.Keys) {
    $ws1.Range("F$r").Value = $exhibitionRows[$r]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($r in $allTypesRows.Keys) {
    $ws4.Range("F$r").Value = $allTypesRows[$r]
}
